$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 297, pushing existing rows 297..337 down to 298..338
$ws.Rows.Item(297).Insert()

# Populate the newly inserted row with the new weekly data point
$ws.Cells.Item(297, 1).Value = 10
$ws.Cells.Item(297, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(297, 3).Value = "La Araucanía"
$ws.Cells.Item(297, 4).Value = 44776
$ws.Cells.Item(297, 5).Value = 9
$ws.Cells.Item(297, 6).Value = 100112017
$ws.Cells.Item(297, 7).Value = "Apio"
$ws.Cells.Item(297, 8).Value = "Americana (o)"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 80
$ws.Cells.Item(297, 11).Value = 12000
$ws.Cells.Item(297, 12).Value = 12000
$ws.Cells.Item(297, 13).Value = 12000
$ws.Cells.Item(297, 14).Value = "$/docena de matas"
$ws.Cells.Item(297, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(297, 16).Value = 2000
$ws.Cells.Item(297, 17).Value = 6
$ws.Cells.Item(297, 18).Value = "Hortaliza"
